$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.274.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.664.53'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5328'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.009'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2638'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06362'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.53'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07824'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.566'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.669.50'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.892.11'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5532'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8210'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.67'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.683'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.38'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.21'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.039'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1227'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.192'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.08'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.486'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05870'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.591'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.277'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.608'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9610'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5798'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01605'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8639'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.048.09'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.07'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.802.63'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.74%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.015'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.23%  '
